$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for new "Std" / "Relative std" block (mirrors the
# existing Min/Max, Q1/Median, Q3/IQR label blocks above it - bold, no
# special number format, same as style index 3 used by those headers).
$ws.Range("D15").Value = "Std"
$ws.Range("E15").Value = "Relative std"
$ws.Range("D15:E15").Font.Bold = $true

# New computed statistics.
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "= (D16 / E4) * 100"

# Leave the active selection on E16 (matches the saved selection in the
# target workbook).
$ws.Range("E16").Select()

$wb.Save()
